$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D and E are treated as text to avoid Excel auto-converting
# numeric-looking strings (e.g. "0.00001050") into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.950.14"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.888.33"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "1.016"
$ws.Range("E4").Value = "  +1.35%  "
$ws.Range("D5").Value = "335.70"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "1.015"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("D7").Value = "0.4701"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "0.3932"
$ws.Range("D9").Value = "46.65"
$ws.Range("E9").Value = "  -3.61%  "
$ws.Range("D10").Value = "0.07990"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "1.899.48"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "5.976"
$ws.Range("D15").Value = "7.158"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").Value = "0.06755"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").Value = "87.70"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "0.00001050"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "17.18"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "27.951.75"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "5.510"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "2.369"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "2.109.70"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").Value = "159.00"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "2.106"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "5.509"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "121.49"
$ws.Range("E31").Value = "  -1.04%  "
$ws.Range("D32").Value = "0.09564"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "0.9653"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").Value = "3.646"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "1.361"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").Value = "0.06135"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "0.02250"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "1.215"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "8.196"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5972"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "0.1902"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "10.35"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.266"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5689"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "12.27"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.399"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.943"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06864"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "113.68"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.069"
$ws.Range("E51").Value = "  -0.61%  "
